# Fixing bug with valence and s value datatypes.
#
# 1) Rename the two sheets/tabs.
# 2) Update the title/"frame of reference" header cells (A1/A2) on each sheet.
# 3) Re-type the "S" column (D) from text (" N") to a real number N, and
#    populate the paired "Valence" column (E) with the matching label -
#    it was previously left blank whenever D held a text value.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- 1) Sheet/tab names -----------------------------------------------
$ws1.Name = "T1 - Test Topic"
$ws2.Name = "T2 - Test2"

# --- 2) Header cells -----------------------------------------------
$ws1.Range("A1").Value = "Test Topic"
$ws1.Range("A2").Value = "Fram"

$ws2.Range("A1").Value = "Test2"
$ws2.Range("A2").Value = "Frame of Reference"

# --- 3) Fix the S (D) / Valence (E) columns -----------------------------------------------
$valenceMap = @{
    1 = "Very Pleasant"
    2 = "Mildly Pleasant"
    3 = "Neutral"
    4 = "Mildly Unpleasant"
    5 = "Very Unpleasant"
}

$sheet1Rows = @{
    10  = 2
    18  = 3
    31  = 3
    41  = 4
    63  = 4
    66  = 4
    67  = 4
    81  = 4
    85  = 4
    106 = 5
    126 = 5
    127 = 5
    130 = 5
    136 = 5
    158 = 5
    167 = 5
    178 = 5
    181 = 5
}

$sheet2Rows = @{
    4   = 1
    5   = 1
    30  = 3
    41  = 3
    51  = 4
    54  = 4
    55  = 4
    56  = 4
    67  = 4
    91  = 4
    92  = 4
    101 = 4
    114 = 4
    128 = 5
    132 = 5
    144 = 5
    145 = 5
    146 = 5
    147 = 5
    173 = 5
    180 = 5
    183 = 5
}

foreach ($row in $sheet1Rows.Keys) {
    $num = $sheet1Rows[$row]
    $ws1.Range("D$row").Value = $num
    $ws1.Range("E$row").Value = $valenceMap[$num]
}

foreach ($row in $sheet2Rows.Keys) {
    $num = $sheet2Rows[$row]
    $ws2.Range("D$row").Value = $num
    $ws2.Range("E$row").Value = $valenceMap[$num]
}
